# Ammend parts list for all-in-one probe:
# Insert a new BOM line above the "Stainless Steel Dowel Pin" row for the
# silicon Ethernet wires that connect to the PCB, pushing the remaining
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts rows 9.. down by one, carrying formulas
# and formatting from the row above, same as Excel's native Insert).
$ws.Rows("9:9").Insert()

# New component line: just a description, no package/value/qty/price yet.
$ws.Range("C9").Value = "Silicon Eithernet Wires to PCB"

# Leave the selection where the user last clicked while editing the list.
$ws.Range("C10").Select()
